$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Lettre Recommandée avec Accusé de Réception"
#       -> "Lettre recommandée avec accusé de réception"
#    Three runs get their casing fixed; every run in the paragraph also
#    becomes italic and gains lang="fr-FR".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Lettre Recommand", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lettre recommand", 2) | Out-Null
$d.Content.Find.Execute("e avec Accus", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e avec accus", 2) | Out-Null
$d.Content.Find.Execute("de R", $true, $false, $false, $false, $false,
                         $true, 1, $false, "de r", 2) | Out-Null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Lettre recommand*") {
        $p.Range.Font.Italic = $true
        $p.Range.LanguageID = "fr-FR"
    }
}

# ---------------------------------------------------------------------------
# 2) "l'expression" (straight apostrophe) -> "l’expression" (curly quote),
#    splitting the run in three (l / ’ / expression...) as is done
#    everywhere else in this letter, matching the source template's
#    convention of giving every curly punctuation mark its own run.
# ---------------------------------------------------------------------------
$needle = $d.Content
$found = $needle.Find.Execute("l'expression", $true, $false, $false, $false, $false,
                               $true, 1, $false, "", 0)
if ($found) {
    $apoStart = $needle.Start + 1
    $apos = $d.Range($apoStart, $apoStart + 1)
    $apos.Text = [string][char]0x2019
    # Re-anchor on the (still 1-char) range and force it to become its own
    # run by toggling a character property - the run boundary persists even
    # after the property is reverted to match its neighbours again.
    $apos = $d.Range($apoStart, $apoStart + 1)
    $apos.Font.Bold = $true
    $apos.Font.Bold = $false
}

Write-Output "done"
